# Add a new "forms" worksheet after the last sheet (calendar), populate it
# with test form data, add mailto hyperlinks on the report-mail column, and
# size the columns to fit the content.

$wb = $excel.ActiveWorkbook

# --- Insert the new sheet after the current last sheet (calendar) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "forms"

# --- Header row ---
$ws.Range("A1").Value = "title"
$ws.Range("B1").Value = "reportMail"
$ws.Range("D1").Value = "weklcomeMessage"
$ws.Range("E1").Value = "confirmationMessage"
$ws.Range("C1").Value = "description"

# --- Column A: title ---
$ws.Range("A2").Value = "Form title  - 1"
$ws.Range("A3").Value = "Form title  - 2"

# --- Column B: reportMail ---
$ws.Range("B2").Value = "tejas.niturkar@gmail.com"
$ws.Range("B3").Value = "pradumna.arts@gmail.com"

# --- Column C: description ---
$ws.Range("C2").Value = "form description - 1 TEST"
$ws.Range("C3").Value = "form description - 2 TEST"

# --- Column D: weklcomeMessage ---
$ws.Range("D2").Value = "Test form welcome message - 1"
$ws.Range("D3").Value = "Test form welcome message - 2"

# --- Column E: confirmationMessage ---
$ws.Range("E2").Value = "Test form confirmation message - 1"
$ws.Range("E3").Value = "Test form confirmation message - 2"

# --- Hyperlinks for the report mail addresses ---
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:tejas.niturkar@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:pradumna.arts@gmail.com")
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"

# --- Column widths (fit to content) ---
$ws.Columns.Item(1).ColumnWidth = 11.8309375
$ws.Columns.Item(2).ColumnWidth = 23.50015625
$ws.Columns.Item(3).ColumnWidth = 32.330625
$ws.Columns.Item(4).ColumnWidth = 28.50015625
$ws.Columns.Item(5).ColumnWidth = 31.830937499999997

# --- Selection on the new sheet ---
$ws.Range("E7").Select()
